# update dayrealmob correct cell in time
# Updates the report date (D1) and the per-hour activity values in column D
# of the "dayforecastonactivity" sheet to the corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date
$ws.Range("D1").Value = "21/03/2023"

# AMM block
$ws.Range("D3").Value = 302.9
$ws.Range("D4").Value = 404.3
$ws.Range("D5").Value = 422.5
$ws.Range("D6").Value = 371.8
$ws.Range("D7").Value = 262.6
$ws.Range("D8").Value = 318.5
$ws.Range("D9").Value = 289.9
$ws.Range("D10").Value = 338
$ws.Range("D11").Value = 315.9
$ws.Range("D12").Value = 343.2
$ws.Range("D13").Value = 256.1
$ws.Range("D14").Value = 166.4

# IPR block
$ws.Range("D18").Value = 9
$ws.Range("D19").Value = 31
$ws.Range("D20").Value = 41
$ws.Range("D21").Value = 30
$ws.Range("D22").Value = 37
$ws.Range("D23").Value = 27
$ws.Range("D24").Value = 33
$ws.Range("D25").Value = 23
$ws.Range("D26").Value = 33
$ws.Range("D27").Value = 23
$ws.Range("D28").Value = 35
$ws.Range("D29").Value = 15
$ws.Range("D30").Value = 21

# MIG block
$ws.Range("D34").Value = 48
$ws.Range("D35").Value = 72

# MOB block
$ws.Range("D36").Value = 89
$ws.Range("D37").Value = 180
$ws.Range("D38").Value = 235
$ws.Range("D39").Value = 229
$ws.Range("D40").Value = 193
$ws.Range("D41").Value = 159
$ws.Range("D42").Value = 180
$ws.Range("D43").Value = 156
$ws.Range("D44").Value = 165
$ws.Range("D45").Value = 175
$ws.Range("D46").Value = 167
$ws.Range("D47").Value = 140
$ws.Range("D48").Value = 83
$ws.Range("D49").Value = 52
$ws.Range("D50").Value = 25
$ws.Range("D51").Value = 18

# MOB PRE block
$ws.Range("D52").Value = 72
$ws.Range("D53").Value = 130
$ws.Range("D54").Value = 162
$ws.Range("D55").Value = 152
$ws.Range("D56").Value = 126
$ws.Range("D57").Value = 106
$ws.Range("D58").Value = 122
$ws.Range("D59").Value = 110
$ws.Range("D60").Value = 126
$ws.Range("D61").Value = 109
$ws.Range("D62").Value = 123
$ws.Range("D63").Value = 90

# MSK block
$ws.Range("D68").Value = 25
$ws.Range("D69").Value = 42
$ws.Range("D70").Value = 43
$ws.Range("D71").Value = 41
$ws.Range("D72").Value = 45
$ws.Range("D73").Value = 40
$ws.Range("D74").Value = 37
$ws.Range("D75").Value = 43
$ws.Range("D76").Value = 40
$ws.Range("D77").Value = 49
$ws.Range("D78").Value = 38
$ws.Range("D79").Value = 31
$ws.Range("D80").Value = 16
$ws.Range("D81").Value = 10
$ws.Range("D82").Value = 4

# NOT block
$ws.Range("D84").Value = 10
$ws.Range("D87").Value = 1
$ws.Range("D88").Value = 2
$ws.Range("D90").Value = 15
$ws.Range("D91").Value = 82

# TEC block
$ws.Range("D92").Value = 100
$ws.Range("D93").Value = 154
$ws.Range("D94").Value = 187
$ws.Range("D95").Value = 173
$ws.Range("D96").Value = 174
$ws.Range("D97").Value = 141
$ws.Range("D98").Value = 169
$ws.Range("D99").Value = 143
$ws.Range("D100").Value = 151
$ws.Range("D101").Value = 157
$ws.Range("D102").Value = 179
$ws.Range("D103").Value = 150
$ws.Range("D104").Value = 107
$ws.Range("D105").Value = 71
$ws.Range("D106").Value = 34
$ws.Range("D107").Value = 19

# TST block
$ws.Range("D109").Value = 24
$ws.Range("D110").Value = 26
$ws.Range("D111").Value = 24
$ws.Range("D112").Value = 21
$ws.Range("D113").Value = 19
$ws.Range("D114").Value = 17
$ws.Range("D115").Value = 20
$ws.Range("D116").Value = 22
$ws.Range("D117").Value = 18
$ws.Range("D118").Value = 25
$ws.Range("D119").Value = 16
$ws.Range("D120").Value = 7

# VIP block
$ws.Range("D125").Value = 1
$ws.Range("D128").Value = 1
$ws.Range("D129").Value = 1
$ws.Range("D131").Value = 1
$ws.Range("D134").Value = 1
$ws.Range("D135").Value = 1

# WLC block
$ws.Range("D140").Value = 7
$ws.Range("D141").Value = 18
$ws.Range("D142").Value = 19
$ws.Range("D143").Value = 24
$ws.Range("D144").Value = 17
$ws.Range("D145").Value = 13
$ws.Range("D146").Value = 17
$ws.Range("D147").Value = 21
$ws.Range("D148").Value = 21
$ws.Range("D149").Value = 37
$ws.Range("D150").Value = 27
$ws.Range("D151").Value = 20
$ws.Range("D152").Value = 11
$ws.Range("D153").Value = 5
$ws.Range("D154").Value = 2
